# Rewrite the "List1" sheet: the ad-hoc exploratory columns collected via
# `Get-ComputerInfo` / `systeminfo` are replaced by a single clean dataset
# (header row + one data row) gathered with a PowerShell script.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Start from a clean sheet - drop everything that's there today.
$ws.Cells.ClearContents()

# Header row
$ws.Range("A1").Value = "WindowsProductName"
$ws.Range("B1").Value = "WindowsVersion"
$ws.Range("C1").Value = "BiosManufacturer"
$ws.Range("D1").Value = "CsDNSHostName"
$ws.Range("E1").Value = "CsDomain"
$ws.Range("F1").Value = "OsTotalVisibleMemorySize"
$ws.Range("G1").Value = "OsArchitecture"
$ws.Range("H1").Value = "Hotfix(s)"

# Data row
$ws.Range("A2").Value = "Windows 10 Home"
$ws.Range("B2").Value = 2009
$ws.Range("C2").Value = "Dell Inc."
$ws.Range("D2").Value = "DESKTOP-SQSHR1A"
$ws.Range("E2").Value = "WORKGROUP"
$ws.Range("F2").Value = 16671872
$ws.Range("G2").Value = "64bitový"
$ws.Range("H2").Value = "KB5012117"

# Match column content with best-fit widths, like Excel does automatically
# when you double-click a column border.
$ws.Columns.Item("A:H").AutoFit()

# Page setup for printing the final sheet.
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# Leave the selection where the author left it before saving.
$ws.Range("E7").Select()
